# conditional_formatting.xlsx — add a second data row with a SUM formula and
# a new conditional-formatting rule that highlights it when > 10.
#
# Fix formula evaluation caching issues (#17854)
#  * Use a common FormulaEvaluator
#  * Force formula evaluation before conditional rule checks
#  * Fixes SHEET-4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data row (row 2 is intentionally left blank, data lands on row 3) ---
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("A3").Formula = "=B3+C3"

# --- new conditional formatting rule on A3: highlight when value > 10 -------
# XlFormatConditionType.xlCellValue = 1, XlFormatConditionOperator.xlGreater = 5
$rngA3 = $ws.Range("A3")
$cfA3 = $rngA3.FormatConditions.Add(1, 5, 10)

# Give it the same "light red fill / dark red text" look already used by the
# existing rule on A1:B1 (font FF9C0006 on fill FFFFC7CE, stored as BGR OLE
# colors for the COM Color properties).
$cfA3.Font.Color = 393372
$cfA3.Interior.Color = 13551615

# The new rule on A3 should be evaluated first (higher priority than the
# pre-existing A1:B1 rule).
$cfA3.SetFirstPriority()

# Move the selection to the newly added cell, like a user would after typing
# the formula.
$ws.Range("A3").Select()
